$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - column headers (bold, matching the existing "title" style used by A3)
$ws.Range("B9").Value = "Number of employees"
$ws.Range("C9").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D9").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B9:D9").Font.Bold = $true

# Row 10 - Micro
$ws.Range("A10").Value = "Micro"
$ws.Range("B10").Value = "1-3"
$ws.Range("C10").Value = "'"
$ws.Range("D10").Value = "'=< FCFA 10,000,000"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Style = "Normal"

# Row 11 - Small
$ws.Range("A11").Value = "Small"
$ws.Range("B11").Value = "4-10"
$ws.Range("C11").Value = "'"
$ws.Range("D11").Value = "FCFA 10,000,000 - 50,000,000"
$ws.Range("C11").Style = "Normal"

# Row 12 - Medium
$ws.Range("A12").Value = "Medium"
$ws.Range("B12").Value = "11-50"
$ws.Range("C12").Value = "'"
$ws.Range("D12").Value = "FCFA 50,000,000 - 500,000,000"
$ws.Range("C12").Style = "Normal"

# Row 13 - Large
$ws.Range("A13").Value = "Large"
$ws.Range("B13").Value = ">50"
$ws.Range("C13").Value = "'"
$ws.Range("D13").Value = "> FCFA 500,000,000"
$ws.Range("C13").Style = "Normal"
